# Update 13C-MFA files (run and result) for SC and IO under WT-batch and chemostats
#
# This script reproduces, via Excel COM interop, the edit captured in the
# OOXML diff: a new flux row ("EX_glc__D_e.f") is inserted into the
# FluxData sheet (shifting every following row down by one), a few of the
# numeric values around the insertion point are updated, and the saved
# view-state (active sheet/tab, zoom, selection) for all three sheets is
# refreshed to match what Excel/LO recorded after the edit.

$wb = $excel.ActiveWorkbook

$msData = $wb.Worksheets.Item("MSData")
$fluxData = $wb.Worksheets.Item("FluxData")
$tracerData = $wb.Worksheets.Item("TracerData")

# ---------------------------------------------------------------------
# FluxData: insert a new row for "EX_glc__D_e.f" right after "BIOMASS.f"
# (row 2), pushing the old row 3 ("EX_c5sugal_e.f") and everything below
# it down by one row.
# ---------------------------------------------------------------------
$fluxData.Rows.Item(3).Insert()

# New row 3: EX_glc__D_e.f
$fluxData.Range("A3").Value = "EX_glc__D_e.f"
$fluxData.Range("B3").Value = 8.49012987008866
$fluxData.Range("C3").Value = 5.63357266159833
$fluxData.Rows.Item(3).RowHeight = 13.8

# Row 2 (BIOMASS.f): only the dilution value changes.
$fluxData.Range("C2").Value = 0.0001

# Row 4 (now EX_c5sugal_e.f, shifted down from row 3): values updated.
$fluxData.Range("B4").Value = 0.005787037037037
$fluxData.Range("C4").Value = 0.026288101469535

# The row that used to be the last one (row 37) is no longer last, so it
# reverts to the regular row height; the new last row (38, shifted down
# from the old row 37 / "DIL_ade_d1.f") takes on the special last-row
# height instead.
$fluxData.Rows.Item(37).RowHeight = 15
$fluxData.Rows.Item(38).RowHeight = 13.8

# ---------------------------------------------------------------------
# View state: zoom every sheet to 100% (was 65%).
# ---------------------------------------------------------------------
$msData.Activate()
[void]($excel.ActiveWindow.Zoom = 100)
$msData.Range("F19").Select()

$tracerData.Activate()
[void]($excel.ActiveWindow.Zoom = 100)
$tracerData.Range("B1").Select()

# FluxData becomes the active/selected sheet (activeTab=1), scrolled back
# to the top with A1 selected.
$fluxData.Activate()
[void]($excel.ActiveWindow.Zoom = 100)
$fluxData.Range("A1").Select()
